$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 becomes "Zusatzbeitrag Krankenversicherung AG-Anteil in Prozent" / 0.99
$ws.Range("A4").Value = "Zusatzbeitrag Krankenversicherung AG-Anteil in Prozent"
$ws.Range("B4").Value = 0.99

# Row 5 becomes "Zusatzbeitrag Krankenversicherung AN-Anteil in Prozent" / 0.99
$ws.Range("A5").Value = "Zusatzbeitrag Krankenversicherung AN-Anteil in Prozent"
$ws.Range("B5").Value = 0.99

# Row 6 becomes "Umlage U1 in Prozent" / 2.3 (previously row 5's content)
$ws.Range("A6").Value = "Umlage U1 in Prozent"
$ws.Range("B6").Value = 2.2999999999999998

# Row 7 becomes "Umlage U2 in Prozent" / 0.44 (previously row 6's content)
$ws.Range("A7").Value = "Umlage U2 in Prozent"
$ws.Range("B7").Value = 0.44

# Row 8 becomes "Insolvenzgeldumlage" / 0.06 (previously row 7's content), now numeric
$ws.Range("A8").Value = "Insolvenzgeldumlage"
$ws.Range("B8").NumberFormat = "0.00"
$ws.Range("B8").Value = 0.06

# New row 9: "Eintragungsdatum" / "01.01.2024" (date kept as text)
$ws.Range("A9").Value = "Eintragungsdatum"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "01.01.2024"

$excel.Goto($ws.Range("B10"))
